# Auto-generated Excel COM-interop script
# Applies the 2026-01-15 12:27:08 scrape update to the "LP1912" / "LP1912-215" / "6203-6173" sheets
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Sheet 1: LP1912
# ---------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("LP1912")
$ws1.Cells.Item(2,1).Value = "Última actualización: 12:27:08"
$ws1.Cells.Item(3,1).Value = "Total filas: 228"

$rows1 = @{
    47 = @("05:49:10", "07:32", "16_SANTA ANA", 103, "LP1912")
    49 = @("05:49:10", "07:32", "11_ETCHEVERRY", 103, "LP1912")
    62 = @("06:43:12", "08:03", "23_HERNANDEZ", 80, "LP1912")
    63 = @("07:19:37", "08:03", "11_ETCHEVERRY", 44, "LP1912")
    105 = @("07:45:49", "09:23", "11_ETCHEVERRY", 98, "LP1912")
    106 = @("07:58:19", "09:23", "17_ROMERO", 85, "LP1912")
    146 = @("10:52:48", "10:56", "16_SANTA ANA", 4, "LP1912")
    147 = @("10:12:35", "10:56", "27_EL RETIRO", 44, "LP1912")
    175 = @("10:12:35", "11:58", "16_P MOR-167 Y 521", 106, "LP1912")
    176 = @("10:52:48", "11:58", "225_GOMEZ", 66, "LP1912")
    190 = @("10:52:48", "12:09", "84_COLONIA URQUIZA-ESC 49", 77, "LP1912")
    191 = @("11:17:08", "12:09", "10_OLMOS", 52, "LP1912")
    201 = @("12:27:08", "12:27", "16_SANTA ANA", 0, "LP1912")
    202 = @("11:45:01", "12:35", "11_ETCHEVERRY", 50, "LP1912")
    203 = @("11:45:01", "12:35", "23_HERNANDEZ", 50, "LP1912")
    204 = @("10:52:48", "12:36", "27_EL RETIRO", 104, "LP1912")
    205 = @("11:17:08", "12:37", "27_EL RETIRO", 80, "LP1912")
    206 = @("10:52:48", "12:38", "17_179 Y 38", 106, "LP1912")
    207 = @("11:17:08", "12:41", "10_OLMOS", 84, "LP1912")
    208 = @("12:27:08", "12:47", "16_SANTA ANA", 20, "LP1912")
    209 = @("10:52:48", "12:48", "11_ETCHEVERRY", 116, "LP1912")
    210 = @("11:17:08", "12:48", "17_ROMERO", 91, "LP1912")
    211 = @("11:17:08", "12:49", "11_ETCHEVERRY", 92, "LP1912")
    212 = @("12:27:08", "12:50", "15_ABASTO", 23, "LP1912")
    213 = @("12:27:08", "12:55", "10_OLMOS", 28, "LP1912")
    214 = @("11:45:01", "13:02", "15_ABASTO", 77, "LP1912")
    215 = @("12:27:08", "13:03", "14_ABASTO", 36, "LP1912")
    216 = @("11:17:08", "13:07", "16_P MOR-SANTA ANA", 110, "LP1912")
    217 = @("12:27:08", "13:08", "10_OLMOS", 41, "LP1912")
    218 = @("11:17:08", "13:14", "215D_EL PATO", 117, "LP1912")
    219 = @("11:17:08", "13:16", "11_ETCHEVERRY", 119, "LP1912")
    220 = @("11:45:01", "13:17", "17_ROMERO", 92, "LP1912")
    221 = @("11:45:01", "13:20", "10_OLMOS", 95, "LP1912")
    222 = @("11:45:01", "13:21", "26_HERNANDEZ", 96, "LP1912")
    223 = @("11:45:01", "13:27", "14_ABASTO", 102, "LP1912")
    224 = @("11:59:06", "13:31", "17_ROMERO", 92, "LP1912")
    225 = @("12:27:08", "13:32", "10_OLMOS", 65, "LP1912")
    226 = @("12:27:08", "13:37", "23_HERNANDEZ", 70, "LP1912")
    227 = @("12:27:08", "13:46", "17_ROMERO", 79, "LP1912")
    228 = @("11:59:06", "13:51", "215A_EL PATO", 112, "LP1912")
    229 = @("11:59:06", "13:56", "225_GOMEZ", 117, "LP1912")
    230 = @("11:59:06", "13:57", "16_P MOR-167 Y 521", 118, "LP1912")
    231 = @("12:27:08", "14:04", "17_ROMERO", 97, "LP1912")
    232 = @("12:27:08", "14:17", "27_EL RETIRO", 110, "LP1912")
    233 = @("12:27:08", "14:20", "215C_EL PATO", 113, "LP1912")
}
foreach ($r in $rows1.Keys) {
    $vals = $rows1[$r]
    $ws1.Cells.Item($r, 1).Value = $vals[0]
    $ws1.Cells.Item($r, 2).Value = $vals[1]
    $ws1.Cells.Item($r, 3).Value = $vals[2]
    $ws1.Cells.Item($r, 4).Value = $vals[3]
    $ws1.Cells.Item($r, 5).Value = $vals[4]
}

# ---------------------------------------------------------------
# Sheet 2: LP1912-215 (append new scraped row 30)
# ---------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("LP1912-215")
$ws2.Cells.Item(2,1).Value = "Última actualización: 12:27:08"
$ws2.Cells.Item(3,1).Value = "Total filas: 25"
$ws2.Cells.Item(30, 1).Value = "12:27:08"
$ws2.Cells.Item(30, 2).Value = "14:20"
$ws2.Cells.Item(30, 3).Value = "215C_EL PATO"
$ws2.Cells.Item(30, 4).Value = 113
$ws2.Cells.Item(30, 5).Value = "LP1912"

# ---------------------------------------------------------------
# Sheet 3: 6203-6173 (append new scraped row 38)
# ---------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("6203-6173")
$ws3.Cells.Item(2,1).Value = "Última actualización: 12:27:08"
$ws3.Cells.Item(3,1).Value = "Total filas: 33"
$ws3.Cells.Item(38, 1).Value = "12:27:08"
$ws3.Cells.Item(38, 2).Value = "14:09"
$ws3.Cells.Item(38, 3).Value = "215A_LA PLATA"
$ws3.Cells.Item(38, 4).Value = 102
$ws3.Cells.Item(38, 5).Value = "L6173"

Write-Output "Update applied: LP1912 rows patched/appended, LP1912-215 row 30 added, 6203-6173 row 38 added."
